$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a couple of existing cells ---------------------------------------
# Row 46: GRADE column was "C3", correct it to "C2".
$ws.Range("F46").Value = "C2"

# Row 199: PUR_INV_ID column had a stray lowercase "i9"; normalize to "I9"
# (matching the rest of the column).
$ws.Range("G199").Value = "I9"

# --- Widen column C slightly ----------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 28.71

# --- Append the new purchase-invoice rows (opening/closing balances) ------
$newRows = @(
    @{ Row=200; A="01.02.21"; B="F-138"; C="M/S.S K M BOARDS&PLYWOODS"; D="7X2.5"; E=1250; F="D1";  G=$null },
    @{ Row=201; A="10.02.21"; B="F-140"; C="M/S.PREMIER EXIM";          D="7X4";   E=7450; F="D1";  G=$null },
    @{ Row=202; A="10.02.21"; B="F-140"; C="M/S.PREMIER EXIM";          D="7X2.5"; E=1000; F="D1";  G=$null },
    @{ Row=203; A="16.02.21"; B="F-142"; C="M/S.FAROOK BOARDS";         D="8X4";   E=1000; F="C3"; G="I9" },
    @{ Row=204; A="16.02.21"; B="F-142"; C="M/S.FAROOK BOARDS";         D="7X4";   E=500;  F="C3"; G="I9" },
    @{ Row=205; A="16.02.21"; B="F-142"; C="M/S.FAROOK BOARDS";         D="6X4";   E=300;  F="C3"; G="I9" },
    @{ Row=206; A="16.02.21"; B="F-142"; C="M/S.FAROOK BOARDS";         D="6X3";   E=300;  F="C3"; G="I9" },
    @{ Row=207; A="16.02.21"; B="F-142"; C="M/S.FAROOK BOARDS";         D="5X3";   E=450;  F="C3"; G="I9" },
    @{ Row=208; A="18.02.21"; B="F-143"; C="M/S.ASHIRVAD PLY";          D="8X4";   E=50;   F="D1";  G=$null }
)

# The "BILL DATE" column holds plain text like "01.02.21" / "10.02.21" in
# this workbook (never a real Excel date). Some of these strings parse as a
# valid date under the default locale, so Excel would silently convert them
# to a date serial when assigned directly. Force the cells to Text format
# first so the literal string is preserved, matching every other row above.
$ws.Range("A200:A208").NumberFormat = "@"

# Fill in column A/B/C first (in the same order the values were first typed
# - Bill Date, then Party Name, then Bill No for the 201/202 block - so new
# entries land in the shared-string table the same way), then backfill the
# remaining columns per row.
$ws.Cells.Item(200, 1).Value = "01.02.21"
$ws.Cells.Item(200, 2).Value = "F-138"
$ws.Cells.Item(200, 3).Value = "M/S.S K M BOARDS&PLYWOODS"

$ws.Cells.Item(201, 1).Value = "10.02.21"
$ws.Cells.Item(202, 1).Value = "10.02.21"
$ws.Cells.Item(201, 3).Value = "M/S.PREMIER EXIM"
$ws.Cells.Item(202, 3).Value = "M/S.PREMIER EXIM"
$ws.Cells.Item(201, 2).Value = "F-140"
$ws.Cells.Item(202, 2).Value = "F-140"

$ws.Cells.Item(203, 1).Value = "16.02.21"
$ws.Cells.Item(204, 1).Value = "16.02.21"
$ws.Cells.Item(205, 1).Value = "16.02.21"
$ws.Cells.Item(206, 1).Value = "16.02.21"
$ws.Cells.Item(207, 1).Value = "16.02.21"
$ws.Cells.Item(203, 2).Value = "F-142"
$ws.Cells.Item(204, 2).Value = "F-142"
$ws.Cells.Item(205, 2).Value = "F-142"
$ws.Cells.Item(206, 2).Value = "F-142"
$ws.Cells.Item(207, 2).Value = "F-142"
$ws.Cells.Item(203, 3).Value = "M/S.FAROOK BOARDS"
$ws.Cells.Item(204, 3).Value = "M/S.FAROOK BOARDS"
$ws.Cells.Item(205, 3).Value = "M/S.FAROOK BOARDS"
$ws.Cells.Item(206, 3).Value = "M/S.FAROOK BOARDS"
$ws.Cells.Item(207, 3).Value = "M/S.FAROOK BOARDS"

$ws.Cells.Item(208, 1).Value = "18.02.21"
$ws.Cells.Item(208, 2).Value = "F-143"
$ws.Cells.Item(208, 3).Value = "M/S.ASHIRVAD PLY"

# Remaining columns (D/E/F/G) all reuse grades/sizes/dates that already
# exist elsewhere in the sheet, so ordering here is not significant.
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    if ($r.G) {
        $ws.Cells.Item($row, 7).Value = $r.G
    }
}

# --- Copy cell formatting (borders/number styles) from matching donor rows -
# Column A uses the "date" style seen on rows 139-199.
$ws.Range("A199").Copy() | Out-Null
$ws.Range("A200:A208").PasteSpecial(-4122) | Out-Null

# Column D uses two border styles seen through the existing data; rows
# 205/206 (6X4 / 6X3) reuse the style from row 190/191, the rest reuse the
# style from row 199.
$ws.Range("D199").Copy() | Out-Null
$ws.Range("D200:D204").PasteSpecial(-4122) | Out-Null
$ws.Range("D207:D208").PasteSpecial(-4122) | Out-Null

$ws.Range("D190").Copy() | Out-Null
$ws.Range("D205:D206").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Update the view / selection so it matches what was on-screen ---------
$ws.Range("H51").Select()
